$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update text content for the quiz/homework lookup response
$ws.Range("E5").Value = "Answer the quiz / homework  /exercises papers"

# Update text content for the schedule update response
$ws.Range("E4").Value = " Update the schedule"

# Update text content (fixing "his or her" -> "his/her")
$ws.Range("B7").Value = "Student wants to check his/her summary of final grade."
$ws.Range("B3").Value = "Student checks his/her schedule"

# Update the view: scroll position and active cell selection
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D4").Select()
